$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Actual Production) values for rows 2-193, taken from the new data export
$bValues = @(351,344,0,343,344,342,344,343,353,348,430,428,431,428,430,525,541,544,552,561,637,537,0,558,725,713,0,705,709,711,710,707,583,573,570,566,412,392,393,374,348,338,327,328,271,239,235,0,238,291,310,436,453,461,490,489,594,611,551,559,603,621,636,649,714,725,759,872,1027,1028,1005,1010,1136,1063,1061,1059,1053,1057,1059,1061,1073,1039,1033,1032,938,900,899,893,646,633,634,625,401,380,386,383,333,329,328,326,328,329,330,328,0,327,328,326,327,0,326,327,326,325,321,311,336,340,0,0,483,433,0,436,580,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($row = 2; $row -le 193; $row++) {
    # Column A timestamps shift forward by exactly 15 days (same time-of-day),
    # matching the new reporting window while keeping full floating point precision.
    $cellA = $ws.Cells.Item($row, 1)
    $oldDate = $cellA.Value()
    $cellA.Value = $oldDate.AddDays(15)

    # Column B gets the refreshed "Actual Production (MW)" readings.
    $ws.Cells.Item($row, 2).Value = $bValues[$row - 2]
}
